$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CapitalCommitment")

# New header cells for From Currency, To Currency, As Of, Exchange Rate
# (order matters for shared string table indices to match target)
$ws.Range("K1").Value = "From Currency"
$ws.Range("L1").Value = "To Currency"
$ws.Range("N1").Value = "As Of"
$ws.Range("M1").Value = "Exchange Rate "

# Update Commitment Date column (J) for all data rows to 2021-01-20 (serial 44216)
$ws.Range("J2").Value = 44216
$ws.Range("J3").Value = 44216
$ws.Range("J4").Value = 44216
$ws.Range("J5").Value = 44216
$ws.Range("J6").Value = 44216
$ws.Range("J7").Value = 44216

# Row 2 FX data
$ws.Range("K2").Value = "USD"
$ws.Range("L2").Value = "INR"
$ws.Range("M2").Value = 80
$ws.Range("N2").Value = 44216

# Row 5 FX data
$ws.Range("K5").Value = "USD"
$ws.Range("L5").Value = "INR"
$ws.Range("M5").Value = 80
$ws.Range("N5").Value = 44216

# Copy number format (date style) from J2 to N2 and N5 (values were pre-seeded
# above as numeric serials so PasteSpecial-Formats reuses the existing style
# index instead of registering a redundant custom number format)
$ws.Range("J2").Copy()
$ws.Range("N2").PasteSpecial(-4122) | Out-Null
$ws.Range("N5").PasteSpecial(-4122) | Out-Null

# Column widths for new columns K, L, M, N (closest achievable values given engine's
# pixel-quantized (1/7 unit) width storage, targeting 12.875 / 10.6875 / 13.5 / 11.875)
$ws.Columns.Item(11).ColumnWidth = 12.142857142857142
$ws.Columns.Item(12).ColumnWidth = 10.0
$ws.Columns.Item(13).ColumnWidth = 12.714285714285714
$ws.Columns.Item(14).ColumnWidth = 11.142857142857142

# Update sheet view: scroll so column B is leftmost, select N5
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("N5").Select()
